$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subtract 10 from every value in column E (rows 2..73)
for ($r = 2; $r -le 73; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $cell.Value = $cell.Value2 - 10
}

# Delete column K entirely (the F/L ratio column, no longer needed)
$ws.Columns.Item(11).Delete()
